$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rescale the count in row 2 (group 0)
$ws.Range("B2").Value = 61

# Copy formatting (bold, border, centered) from the existing labeled cell A3
# onto the new label cells A4:A5 that the extra rows introduce.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-write the group/time labels and counts for the expanded table:
# group 3 now occupies row 3, the former group-1 row becomes row 4 with its
# rescaled count, and a brand new group 2 row is appended as row 5.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 46

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 31

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 26
